$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 439, shifting the existing rows 439-468 down to 441-470.
$ws.Range("A439:A440").EntireRow.Insert()

# New row 439 data
$ws.Range("A439").Value = 8
$ws.Range("B439").Value = "Terminal La Palmera de La Serena"
$ws.Range("C439").Value = "Coquimbo"
$ws.Range("D439").Value = 45223
$ws.Range("E439").Value = 4
$ws.Range("F439").Value = 100112031
$ws.Range("G439").Value = "Poroto verde"
$ws.Range("H439").Value = "Magnum"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 440
$ws.Range("K439").Value = 32000
$ws.Range("L439").Value = 33000
$ws.Range("M439").Value = 32500
$ws.Range("N439").Value = "`$/caja 25 kilos"
$ws.Range("O439").Value = "Provincia de Limarí"
$ws.Range("P439").Value = 1300
$ws.Range("Q439").Value = 25
$ws.Range("R439").Value = "Hortaliza"

# New row 440 data
$ws.Range("A440").Value = 8
$ws.Range("B440").Value = "Terminal La Palmera de La Serena"
$ws.Range("C440").Value = "Coquimbo"
$ws.Range("D440").Value = 45223
$ws.Range("E440").Value = 4
$ws.Range("F440").Value = 100112031
$ws.Range("G440").Value = "Poroto verde"
$ws.Range("H440").Value = "Magnum"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 400
$ws.Range("K440").Value = 25000
$ws.Range("L440").Value = 26000
$ws.Range("M440").Value = 25500
$ws.Range("N440").Value = "`$/malla 25 kilos"
$ws.Range("O440").Value = "Perú"
$ws.Range("P440").Value = 1020
$ws.Range("Q440").Value = 25
$ws.Range("R440").Value = "Hortaliza"
